# Adding new iAuthor testcases
# Replace the existing test-result rows with a single new test case row,
# shrinking the used range from A1:C6 down to A1:C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unneeded data rows (3-6), leaving only the header (row 1)
# and a single data row (row 2).
$ws.Rows("3:6").Delete()

# Overwrite the remaining data row with the new iAuthor test case.
$ws.Range("A2").Value = "iAU_TC_ID_48"
$ws.Range("B2").Value = "@RegressionA Validation of Create Question (Type-B) - Negative Scenario."
$ws.Range("C2").Value = "timedOut"
